# New weekly "Crespo record" price observation for
# Terminal Hortofrutícola Agro Chillán - Repollo.
#
# The diff shows a row inserted at sheet row 166 (old rows 166-177 shift
# down to 167-178, dimension grows from A1:R177 to A1:R178). Insert a
# whole row there, then populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(166).Insert()

$ws.Cells.Item(166, 1).Value  = 7
$ws.Cells.Item(166, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(166, 3).Value  = "Ñuble"
$ws.Cells.Item(166, 4).Value  = 44585
$ws.Cells.Item(166, 5).Value  = 16
$ws.Cells.Item(166, 6).Value  = 100112006
$ws.Cells.Item(166, 7).Value  = "Repollo"
$ws.Cells.Item(166, 8).Value  = "Crespo record"
$ws.Cells.Item(166, 9).Value  = "Primera"
$ws.Cells.Item(166, 10).Value = 200
$ws.Cells.Item(166, 11).Value = 650
$ws.Cells.Item(166, 12).Value = 700
$ws.Cells.Item(166, 13).Value = 675
$ws.Cells.Item(166, 14).Value = "`$/unidad"
$ws.Cells.Item(166, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(166, 16).Value = 675
$ws.Cells.Item(166, 17).Value = 1
$ws.Cells.Item(166, 18).Value = "Hortaliza"
